# Automatische test-sync: 2025-07-27 19:34:50
# Appends the newest "Retour / Terugbetaling" test-mail log entry (row 13)
# to the Logs sheet, rolls the per-category tally on Dashboard (row 6),
# and extends the Dashboard bar chart's category/value series to include it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet: append row 13 with the new mail-log record
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A13").Value = "Mijn retour is nog steeds niet verwerkt."
$logs.Range("B13").Value = "mailmind.test@zohomail.eu"
$logs.Range("C13").Value = "Testmail #11: Mijn retour is nog steeds niet verwerkt."
$logs.Range("D13").Value = "Retour / Terugbetaling"
$logs.Range("E13").Value = "Geachte klant,`nDank u wel voor uw e-mail. Om u beter van dienst te kunnen zijn, ontvangen wij graag wat meer informatie om uw retourzending te kunnen traceren. Kunt u ons uw ordernummer en de datum van retournering verstrekken? Op die manier kunnen we uw zaak nauwkeurig onderzoeken en u van dienst zijn.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F13").Value = "2025-07-27 19:34:37"
$logs.Range("G13").Value = "Ja"
$logs.Range("H13").Value = "Nee"
$logs.Range("I13").Value = "Ja"
$logs.Range("J13").Value = "Nee"

# Conditional formatting ranges D:J grow from row 12 to row 13 to cover the new
# row. Reuse the existing rules (same type/operator/dxf/priority) and just widen
# the sqref each applies to, rather than deleting and re-adding them.
$cfCols = @("D", "G", "H", "I", "J")
foreach ($col in $cfCols) {
    $oldRange = $logs.Range($col + "2:" + $col + "12")
    $newRange = $logs.Range($col + "2:" + $col + "13")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 2) Dashboard sheet: append the per-category tally for the new category
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A6").Value = "Retour / Terugbetaling"
$dash.Range("B6").Value = 1

# ---------------------------------------------------------------------
# 3) Dashboard chart: extend category/value series to the new row
# ---------------------------------------------------------------------
$chart = $dash.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$6,'Dashboard'!`$B`$2:`$B`$6,1)"
